# Apply the edits described by the commit:
#  1. Insert a new column before "Note" (column DL) carrying the header
#     "DemonstrationProjectIdentifier" - this shifts every column from the
#     old DL onward one position to the right (DL->DM, ... MN->MO).
#  2. Update the Id column (A2:A5) value to the new id.
#  3. Clear the SubscriberRelationshipType value in AT2 (was "SELF").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert the new column and set its header.
$ws.Columns("DL:DL").Insert()
$ws.Range("DL1").Value = "DemonstrationProjectIdentifier"

# 2. Update the Id values for the four data rows.
$ws.Range("A2").Value = "6901488a7e79911955eafdaa"
$ws.Range("A3").Value = "6901488a7e79911955eafdaa"
$ws.Range("A4").Value = "6901488a7e79911955eafdaa"
$ws.Range("A5").Value = "6901488a7e79911955eafdaa"

# 3. Clear the SubscriberRelationshipType value on row 2.
$ws.Range("AT2").ClearContents()
